$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 2002年-2009年 rows (rows 2-9), shifting the remaining rows
# (2011年-2019年, previously rows 10-18) up to rows 2-10.
$ws.Range("A2:G9").Delete(-4162)

# Add the new 2021年 row at row 11, copying the existing row's formatting
# (bold/centered/bordered year cell) so no new cell styles are minted.
$ws.Range("A10").Copy($ws.Range("A11"))

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 1232677
$ws.Range("C11").Value = 607439
$ws.Range("D11").Value = 30145
$ws.Range("E11").Value = 39531
$ws.Range("F11").Value = 625238
$ws.Range("G11").Value = 9386
